$d = $word.ActiveDocument

# --- 1. Code block text: rename metric labels passed to modify_header() ---
# These quoted strings are unique in the document, so a whole-document Find
# (narrowed to just the match) is safe and preserves the surrounding
# syntax-highlighted run structure.
$rng = $d.Content
$rng.Find.Execute('"**T-Statistic**"', $true, $false, $false, $false, $false, $true, 0, $false, "", 0) | Out-Null
$rng.Text = '"**t-statistic**"'

$rng = $d.Content
$rng.Find.Execute('"**P-Value**"', $true, $false, $false, $false, $false, $true, 0, $false, "", 0) | Out-Null
$rng.Text = '"**p-value**"'

# --- 2. Regression table: update header labels and all changed values ---
$t = $d.Tables.Item(1)

# Header row (now unambiguous since the quoted code-block text was already changed)
$t.Cell(1,3).Range.Text = "t-statistic"
$t.Cell(1,5).Range.Text = "p-value"

# Row 2: depression_mean
$t.Cell(2,3).Range.Text = "39.1"
$t.Cell(2,4).Range.Text = "0.78, 0.87"

# Row 3: neuroticism_diff
$t.Cell(3,3).Range.Text = "4.73"

# Row 4: neuroticism_mean
$t.Cell(4,3).Range.Text = "-1.40"
$t.Cell(4,5).Range.Text = "0.2"

# Row 5: sex_1
$t.Cell(5,2).Range.Text = "0.03"
$t.Cell(5,3).Range.Text = "1.05"
$t.Cell(5,4).Range.Text = "-0.03, 0.09"
$t.Cell(5,5).Range.Text = "0.3"

# Row 6: race_1
$t.Cell(6,3).Range.Text = "-3.99"

# Row 7: sex_2
$t.Cell(7,3).Range.Text = "-1.87"
$t.Cell(7,4).Range.Text = "-0.11, 0.00"
$t.Cell(7,5).Range.Text = "0.061"

# Row 9: model fit statistics note -- this row's cell is horizontally merged
# (gridSpan=5), and merged-cell ranges insert rather than overwrite, so use
# the same narrowed whole-document Find/replace trick used above.
$rng = $d.Content
$rng.Find.Execute("0.670", $true, $false, $false, $false, $false, $true, 0, $false, "", 0) | Out-Null
$rng.Text = "0.671"

$rng = $d.Content
$rng.Find.Execute("4,179", $true, $false, $false, $false, $false, $true, 0, $false, "", 0) | Out-Null
$rng.Text = "4,180"

$rng = $d.Content
$rng.Find.Execute("4,224", $true, $false, $false, $false, $false, $true, 0, $false, "", 0) | Out-Null
$rng.Text = "4,225"

$d.Save()
